$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B62").Value = 96.90000000000001
$ws.Range("B63").Value = 99.3
$ws.Range("B64").Value = 100.06
$ws.Range("B68").Value = 99.28
$ws.Range("B75").Value = 100.6
$ws.Range("B80").Value = 101.06
$ws.Range("B81").Value = 103.57
$ws.Range("B83").Value = 94.72
$ws.Range("B87").Value = 106.02
